# Updated Excel with data from Adafruit IO
# Append one new reading row (row 94) to the sheet, matching the format
# of the existing Adafruit IO export rows (Timestamp, Feed Key, Value,
# Latitude, Longitude, Elevation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 94

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# "Value" column holds numeric-looking text (e.g. "25") in this sheet, so
# force it to stay text (matching the other rows) instead of being
# auto-converted to a number, then drop back to the default style so no
# stray number-format is left behind.
$valueCell = $ws.Cells.Item($row, 3)
$valueCell.Value = "'25"
$valueCell.Style = "Normal"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
